$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$rows = $tbl.Rows.Count
$cols = $tbl.Columns.Count

$borderColor = 13882323   # RGB(0xD3,0xD3,0xD3) -> D3D3D3
$headerFill  = 12180223   # RGB(0xFF,0xDA,0xB9) -> FFDAB9

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $tbl.Cell($r, $c)

        # Add a thin single light-gray border on every side of every cell
        # (top/left/bottom/right only -- leave inside borders untouched).
        foreach ($side in -1, -2, -3, -4) {
            $cell.Borders.Item($side).LineStyle = 1
            $cell.Borders.Item($side).LineWidth = 2
            $cell.Borders.Item($side).Color = $borderColor
        }

        if ($r -eq 1) {
            # Header row: change the fill color.
            $cell.Shading.BackgroundPatternColor = $headerFill
        } else {
            if ($c -eq 2) {
                # Data rows, second column: right -> center alignment.
                $cell.Range.ParagraphFormat.Alignment = 1
            }
        }
    }
}

Write-Host "Table formatting updated"
